$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "kelas" column header
$ws.Range("D2").Value = "kelas"

# New "kelas" column (D): first 10 students are XRPL1, remaining 9 are XRPL2
$ws.Range("D3:D12").Value = "XRPL1"
$ws.Range("D13:D21").Value = "XRPL2"

# jurusan_kelas (column C) now holds the full major name for every student
$ws.Range("C3:C21").Value = "REKAYASA PERANGKAT LUNAK"

$ws.Range("B7").Select() | Out-Null
